$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("studyAmendments")

# Insert three new columns before column A (shifts old A:G -> D:J),
# making room for new "name" / "label" / "description" columns.
$ws.Range("A1:C1").EntireColumn.Insert()

# New header row (A1:C1) - copy formatting from an existing header cell
# (now at E1, style matches the other header cells) then left-align it.
$ws.Range("E1").Copy()
$ws.Range("A1:C1").PasteSpecial(-4122)
$ws.Range("A1:C1").HorizontalAlignment = -4131

$ws.Range("A1").Value = "name"
$ws.Range("B1").Value = "label"
$ws.Range("C1").Value = "description"

# New amendment name column values
$ws.Range("A2").Value = "AMEND_1"
$ws.Range("A3").Value = "AMEND_2"
$ws.Range("A4").Value = "AMEND_3"
$ws.Range("A5").Value = "AMEND_4"

# Match left-alignment styling used for the new name/label/description cells
$ws.Range("A2:C2").HorizontalAlignment = -4131
$ws.Range("A3:A5").HorizontalAlignment = -4131

# Add phonetic info (matches workbook change) and make this the active sheet/cell
$ws.PhoneticConversion = 2

$ws.Activate()
$ws.Range("C8").Select()
